# make setup save to excel to keep participants blind to condition
# Adds a new participant row (row 11) to Sheet1 recording a test run
# where setup/quit state is now persisted to the spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 11

$ws.Cells.Item($row, 1).Value = "studyname_20201218_04"   # A file
$ws.Cells.Item($row, 2).Value = 10                          # B participant
$ws.Cells.Item($row, 3).Value = 44183                        # C test_date
# Reuse the date number-format already used by the cell above it
# (copy/paste-format keeps the same style record instead of minting a new one)
$ws.Cells.Item($row - 1, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 4).Value = "MZ"                          # D researcher
$ws.Cells.Item($row, 5).Value = "Zoom"                        # E location
$ws.Cells.Item($row, 10).Value = "condition2"                 # J condition
$ws.Cells.Item($row, 11).Value = "left"                       # K counterbalance
$ws.Cells.Item($row, 17).Value = "no"                         # Q in_progress
$ws.Cells.Item($row, 18).Value = "Yes"                        # R measure1
$ws.Cells.Item($row, 19).Value = "blueberries"                # S measure2
$ws.Cells.Item($row, 20).Value = "blue"                       # T measure3
$ws.Cells.Item($row, 21).Value = "right"                      # U measure4
$ws.Cells.Item($row, 22).Value = "right"                      # V measure5
$ws.Cells.Item($row, 23).Value = "clicked"                    # W measure6_button1
$ws.Cells.Item($row, 24).Value = "clicked"                    # X measure6_button2
$ws.Cells.Item($row, 27).Value = "clicked"                    # AA measure6_button5
$ws.Cells.Item($row, 28).Value = "clicked"                    # AB measure6_button6
$ws.Cells.Item($row, 29).Value = "easy"                       # AC measure7_1
$ws.Cells.Item($row, 30).Value = "very easy"                  # AD measure7_2
$ws.Cells.Item($row, 31).Value = "test quit, resume after setup" # AE measure8
$ws.Cells.Item($row, 32).Value = 2                            # AF measure9_highPerf_anchor
$ws.Cells.Item($row, 33).Value = 4                            # AG measure9_lowPerf_anchor
$ws.Cells.Item($row, 34).Value = 3                            # AH measure10_highPerf_anchor
$ws.Cells.Item($row, 35).Value = 3                            # AI measure10_lowPerf_anchor
$ws.Cells.Item($row, 36).Value = 3                            # AJ consent
$ws.Cells.Item($row, 37).Value = "chromebook"                 # AK device
$ws.Cells.Item($row, 38).Value = 4                            # AL video
$ws.Cells.Item($row, 39).Value = 2                            # AM sound
$ws.Cells.Item($row, 40).Value = 3                            # AN fun
